# Møller Innkjøp stasjoner.xlsx -- "Add files via upload"
# Updates the "Viking" sheet: renames the Oslo station, appends six new
# Viking Kontroll stations (with their coordinates as text), adds a
# formatted-but-empty row 10 cell, resizes the data columns, and leaves the
# selection on B13 (mirrors what Excel's own autosave/upload snapshot shows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Viking")

function Set-TextCell($cell, [string]$text) {
    # Force the cell to stay text (so numeric-looking lat/long strings such
    # as "59.92815" are not silently coerced into floating point numbers),
    # then drop back to the workbook's default "Normal" style so no stray
    # numeric-format style sticks to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Row 2: station renamed from "Vollaveien" to "Viking Kontroll Alna" ---
$ws.Range("A2").Value = "Viking Kontroll Alna"

# --- New rows 3-8: additional Viking Kontroll stations ---
$newRows = @(
    @("Viking Kontroll Fredrikstad", "Stabburveien 6",          "Fredrikstad",    1617, "59.2385261",        "10.9647012"),
    @("Viking Kontroll Gjøvik",      "Bryggevegen 9",            "Gjøvik",         2821, "60.795694",         "10.6992484"),
    @("Viking Kontroll Hamar",       "Halsetsvea 38",            "Ingeberg",       2323, "60.8398567",        "11.0937347"),
    @("Viking Kontroll Kristiansand","Skibåsen 33 A",            "Kristiansand S", 4636, "58.1782423",        "8.1207413"),
    @("Viking Kontroll Trondheim",   "Magnus Lagabøters veg 8",  "Trondheim",      7047, "63.4242057800293",  "10.4752950668335"),
    @("Viking Kontroll Alta",        "Knottveien 1",             "Alta",           9514, "69.9665723",        "23.3634596")
)

$r = 3
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    Set-TextCell $ws.Range("E$r") $row[4]
    Set-TextCell $ws.Range("F$r") $row[5]
    $r++
}

# --- Row 10: leftover formatted (but empty) cell, as seen in the upload ---
$a10 = $ws.Cells.Item(10, 1)
$a10.Font.Color = 0
$a10.Font.Name = "Aptos"

# --- Column widths for the widened station/address/coordinate columns ---
# (tuned so the engine's pixel-rounded stored width lands on the same value
# Excel's own AutoFit produced: 23.1640625 / 21.5 / 12.5 / 12 / 16.83203125)
$ws.Columns.Item(1).ColumnWidth = 22.333333333333332
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 16.0
$ws.Columns.Item(6).ColumnWidth = 16.0

# --- Leave the selection where the author's last save left it ---
$ws.Range("B13").Select()
